# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") values are forced to text with a leading apostrophe so
# Excel doesn't reinterpret dotted/decimal-looking strings as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.974.60"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "'1.645.17"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'217.68"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "'0.5231"
$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'0.2615"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").Value = "'0.06267"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("D10").Value = "'20.43"
$ws.Range("E10").Value = "  -3.05%  "

$ws.Range("D11").Value = "'0.07720"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "'4.452"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "'1.648.53"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "'0.5433"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "'0.0₅8063"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("D16").Value = "'64.62"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "'25.998.54"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "'4.540"
$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").Value = "'191.18"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").Value = "'5.974"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "'139.20"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "'0.1234"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'7.249"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'16.14"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").Value = "'1.426"
$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").Value = "'0.05923"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("E32").Value = "  -3.25%  "

$ws.Range("D33").Value = "'1.516"
$ws.Range("E33").Value = "  -8.02%  "

$ws.Range("D34").Value = "'2.417"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").Value = "'0.9416"
$ws.Range("E35").Value = "  -4.05%  "

$ws.Range("D36").Value = "'2.748"
$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").Value = "'0.5697"
$ws.Range("E37").Value = "  -3.81%  "

$ws.Range("D38").Value = "'0.01601"
$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("D39").Value = "'5.842"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("D40").Value = "'0.8463"
$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("D43").Value = "'1.000.31"
$ws.Range("E43").Value = "  -3.63%  "

$ws.Range("D44").Value = "'1.787.81"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'56.50"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.0₈106"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").Value = "'0.4290"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").Value = "'1.476"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").Value = "'0.05149"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").Value = "'7.839"
$ws.Range("E51").Value = "  -3.35%  "
